$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing values for row 46 (2025-02-12 abs_activity)
$ws.Range("D46").Value = 10
$ws.Range("F46").Value = 20

# Fix existing values for row 47 (2025-02-12 rel_activity)
$ws.Range("D47").Value = 5.464750464750465
$ws.Range("F47").Value = 5.464750464750465

# Add new rows 50-53 for 2025-02-13
$ws.Range("A50:A53").NumberFormat = "@"

$ws.Range("A50").Value = "2025-02-13"
$ws.Range("B50").Value = "abs_activity"
$ws.Range("C50").Value = 9.912783085418916
$ws.Range("D50").Value = 8.086492844885454
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 17.99927593030437

$ws.Range("A51").Value = "2025-02-13"
$ws.Range("B51").Value = "rel_activity"
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0

$ws.Range("A52").Value = "2025-02-13"
$ws.Range("B52").Value = "abs_sleep"
$ws.Range("C52").Value = 10
$ws.Range("D52").Value = 10
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 20

$ws.Range("A53").Value = "2025-02-13"
$ws.Range("B53").Value = "rel_sleep"
$ws.Range("C53").Value = 10
$ws.Range("D53").Value = 8.704137766084667
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 18.70413776608467

$ws.Range("A50:A53").ClearFormats()
